$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.598.52'
$ws.Range('E2').Value = '  -2.73%  '

$ws.Range('D3').Value = '3.005.99'
$ws.Range('E3').Value = '  -2.39%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '546.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.03%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.77'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.64%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.09%  '

$ws.Range('D8').Value = '3.000.85'
$ws.Range('E8').Value = '  -2.34%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.490'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.02%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.01'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.04%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.145'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -7.88%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.445'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.37%  '

$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '34.09'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.94%  '

$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000218'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.67%  '

$ws.Range('D15').Value = '3.486.34'
$ws.Range('E15').Value = '  -2.59%  '

$ws.Range('D16').Value = '61.676.18'
$ws.Range('E16').Value = '  -2.64%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.110'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.70%  '

$ws.Range('D18').Value = '3.000.56'
$ws.Range('E18').Value = '  -2.57%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.62'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.50%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '483.54'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.00%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.19'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.71%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.665'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.15%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.15%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.64'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.78%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.02'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.49%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.02%  '

$ws.Range('E27').Value = '  -0.99%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.65'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.06%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.17%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.91'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.90%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.59'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.37%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.12'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.36%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.34'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.50%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.58'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.52%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '54.73'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.15%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.86'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.46%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '440.65'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -9.74%  '

$ws.Range('D38').Value = '3.134.65'
$ws.Range('E38').Value = '  -3.66%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0793'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.14%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0383'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.79%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.116'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.18%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.07'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.95%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.39'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.72%  '

$ws.Range('E44').Value = '  -0.01%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '25.97'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.71%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.242'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.68%  '

$ws.Range('E47').Value = '  -0.98%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.94'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.52%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '114.95'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.94%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.29'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.40%  '

$ws.Range('D51').Value = '0.0₃0483'
$ws.Range('E51').Value = '  -8.33%  '
